$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 179.73685
$ws.Range("J33").Value = 87
$ws.Range("L33").Value = 87
$ws.Range("N33").Value = -545

# Row 99
$ws.Range("H99").Value = 259.8
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# Row 137
$ws.Range("H137").Value = 1156.9
$ws.Range("I137").Value = 1054.5714
$ws.Range("K137").Value = 3163.7142
$ws.Range("M137").Value = -613.7142000000003

# Row 138
$ws.Range("H138").Value = 6045.9
$ws.Range("I138").Value = 2698.6365
$ws.Range("J138").Value = 7315.552
$ws.Range("K138").Value = 8095.9095
$ws.Range("L138").Value = 21946.656
$ws.Range("M138").Value = -2955.9095
$ws.Range("N138").Value = -32226.656


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2480.6584
$ws.Range("I32").Value = 1297.9166
$ws.Range("K32").Value = 1297.9166
$ws.Range("M32").Value = -1010.9166

# Row 61
$ws.Range("H61").Value = 4312.2085
$ws.Range("I61").Value = 4096.5713
$ws.Range("K61").Value = 4096.5713
$ws.Range("M61").Value = -3884.5713

# Row 110
$ws.Range("H110").Value = 1552
$ws.Range("I110").Value = 1552
$ws.Range("K110").Value = 1552
$ws.Range("M110").Value = 493

# Row 132
$ws.Range("H132").Value = 9515.304
$ws.Range("I132").Value = 9784.796
$ws.Range("J132").Value = 7628.857
$ws.Range("K132").Value = 29354.388
$ws.Range("L132").Value = 22886.571
$ws.Range("M132").Value = -26824.388
$ws.Range("N132").Value = -27946.571

# Row 136
$ws.Range("H136").Value = 4312.2085
$ws.Range("I136").Value = 4096.5713
$ws.Range("K136").Value = 12289.7139
$ws.Range("M136").Value = -9739.713899999999


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 3031.2104
$ws.Range("I94").Value = 1298.625
$ws.Range("K94").Value = 1298.625
$ws.Range("M94").Value = -847.625

# Row 108
$ws.Range("H108").Value = 61997.5
$ws.Range("J108").Value = 61997.5
$ws.Range("L108").Value = 61997.5
$ws.Range("N108").Value = -69677.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 225.22223
$ws.Range("I22").Value = 244.66667
$ws.Range("J22").Value = 128
$ws.Range("K22").Value = 244.66667
$ws.Range("L22").Value = 128
$ws.Range("M22").Value = 105.33333
$ws.Range("N22").Value = -828

# Row 62
$ws.Range("H62").Value = 3824.1738
$ws.Range("I62").Value = 3629.3157
$ws.Range("J62").Value = 4749.75
$ws.Range("K62").Value = 3629.3157
$ws.Range("L62").Value = 4749.75
$ws.Range("M62").Value = -3005.3157
$ws.Range("N62").Value = -5997.75

# Row 65
$ws.Range("H65").Value = 3824.1738
$ws.Range("I65").Value = 3629.3157
$ws.Range("J65").Value = 4749.75
$ws.Range("K65").Value = 18146.5785
$ws.Range("L65").Value = 23748.75
$ws.Range("M65").Value = -15026.5785
$ws.Range("N65").Value = -29988.75

# Row 99
$ws.Range("H99").Value = 1170.925
$ws.Range("I99").Value = 1174.9744
$ws.Range("K99").Value = 1174.9744
$ws.Range("M99").Value = 323.0255999999999

# Row 105
$ws.Range("H105").Value = 1338.5714
$ws.Range("I105").Value = 1384.5
$ws.Range("K105").Value = 1384.5
$ws.Range("M105").Value = 362.5

# Row 126
$ws.Range("H126").Value = 1170.925
$ws.Range("I126").Value = 1174.9744
$ws.Range("K126").Value = 3524.9232
$ws.Range("M126").Value = -1054.9232

# Row 132
$ws.Range("H132").Value = 2342.9443
$ws.Range("I132").Value = 2134.4814
$ws.Range("K132").Value = 6403.4442
$ws.Range("M132").Value = -3873.4442

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1073.2307
$ws.Range("I122").Value = 781.2857
$ws.Range("J122").Value = 1413.8334
$ws.Range("K122").Value = 7031.571300000001
$ws.Range("L122").Value = 12724.5006
$ws.Range("M122").Value = -4581.571300000001
$ws.Range("N122").Value = -17624.5006

# Row 128
$ws.Range("H128").Value = 789661.7
$ws.Range("I128").Value = 789661.7
$ws.Range("K128").Value = 2368985.1
$ws.Range("M128").Value = -2364005.1


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 48545.168
$ws.Range("I70").Value = 90845
$ws.Range("K70").Value = 90845
$ws.Range("M70").Value = -90575

# Row 73
$ws.Range("H73").Value = 48545.168
$ws.Range("I73").Value = 90845
$ws.Range("K73").Value = 90845
$ws.Range("M73").Value = -89909

# Row 80
$ws.Range("H80").Value = 6063.5
$ws.Range("I80").Value = 6063.5
$ws.Range("K80").Value = 6063.5
$ws.Range("M80").Value = -5065.5

# Row 83
$ws.Range("H83").Value = 6063.5
$ws.Range("I83").Value = 6063.5
$ws.Range("K83").Value = 30317.5
$ws.Range("M83").Value = -25325.5

# Row 97
$ws.Range("H97").Value = 906.3333
$ws.Range("I97").Value = 906.3333
$ws.Range("K97").Value = 906.3333
$ws.Range("M97").Value = -410.3333

# Row 101
$ws.Range("H101").Value = 20157
$ws.Range("J101").Value = 20157
$ws.Range("L101").Value = 20157
$ws.Range("N101").Value = -26647

# Row 105
$ws.Range("H105").Value = 14487.5
$ws.Range("J105").Value = 14487.5
$ws.Range("L105").Value = 14487.5
$ws.Range("N105").Value = -21475.5


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1913.8928
$ws.Range("I16").Value = 1180
$ws.Range("J16").Value = 3463.2222
$ws.Range("K16").Value = 1180
$ws.Range("L16").Value = 3463.2222
$ws.Range("M16").Value = -1010
$ws.Range("N16").Value = -3803.2222

# Row 124
$ws.Range("H124").Value = 46686
$ws.Range("J124").Value = 46686
$ws.Range("L124").Value = 46686
$ws.Range("N124").Value = -56506


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4217.9375
$ws.Range("I81").Value = 2398.1538
$ws.Range("K81").Value = 4796.3076
$ws.Range("M81").Value = -3735.3076

# Row 84
$ws.Range("H84").Value = 4217.9375
$ws.Range("I84").Value = 2398.1538
$ws.Range("K84").Value = 23981.538
$ws.Range("M84").Value = -18677.538

# Row 113
$ws.Range("H113").Value = 2175.9
$ws.Range("I113").Value = 2228.7778
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 6686.3334
$ws.Range("L113").Value = 5100
$ws.Range("M113").Value = -4516.3334
$ws.Range("N113").Value = -9440

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 122
$ws.Range("H122").Value = 4606.96
$ws.Range("I122").Value = 4938.778
$ws.Range("K122").Value = 14816.334
$ws.Range("M122").Value = -12366.334

# Row 132
$ws.Range("H132").Value = 2577.0625
$ws.Range("I132").Value = 1889.3636
$ws.Range("K132").Value = 5668.0908
$ws.Range("M132").Value = -3138.0908

# Row 136
$ws.Range("H136").Value = 5615.15
$ws.Range("I136").Value = 2120.2
$ws.Range("K136").Value = 6360.599999999999
$ws.Range("M136").Value = -3810.599999999999

